$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 74, pushing the existing
# rows 74-76 down to 75-77 (formatting/styles are carried along
# automatically, matching the d="2" date style used by column D).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly record
# (same "template" as the neighbouring rows, but with its own
# date / volume / price / origin data).
$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C74").Value = "Los Lagos"
$ws.Range("D74").Value = 44448
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100108
$ws.Range("H74").Value = "Tropicales y subtropicales"
$ws.Range("I74").Value = 100108002
$ws.Range("J74").Value = "Mango"
$ws.Range("K74").Value = "Sin especificar"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 80
$ws.Range("N74").Value = 12000
$ws.Range("O74").Value = 12000
$ws.Range("P74").Value = 12000
$ws.Range("Q74").Value = "$/bandeja 4 kilos"
$ws.Range("R74").Value = "Brasil"
$ws.Range("S74").Value = 3000
$ws.Range("T74").Value = 4
